{"js": "// Issue 3 (\"HTTP error 404 appears after deletion of Contact\") gains a new\n// step in its \"Steps to reproduce\" list: \"Click \"Yes\" on pop-up window\",\n// inserted right after \"On the Contact Info page click \"Delete\" link\" and\n// before \"Open Dev tools/Network\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"On the Contact Info page click \\u201CDelete\\u201D link\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find the 'On the Contact Info page click \\u201CDelete\\u201D link' step.\");\n}\n\n// insertParagraph(..., \"After\") inherits the anchor's paragraph formatting\n// (the \"ListParagraph\" style + numPr numbering), matching how the sibling\n// steps in this list are structured.\nanchor.insertParagraph(\"Click \\u201CYes\\u201D on pop-up window\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Issue 3 (\"HTTP error 404 appears after deletion of Contact\") gains a new\n# step in its \"Steps to reproduce\" list: Click \"Yes\" on pop-up window,\n# inserted right after \"On the Contact Info page click \"Delete\" link\" and\n# before \"Open Dev tools/Network\".\n\n$d = $word.ActiveDocument\n\n$quoteOpen = [char]0x201C\n$quoteClose = [char]0x201D\n$anchorText = \"On the Contact Info page click \" + $quoteOpen + \"Delete\" + $quoteClose + \" link\"\n$newStepText = \"Click \" + $quoteOpen + \"Yes\" + $quoteClose + \" on pop-up window\"\n\n# Locate the anchor paragraph with Find (confirms the text exists in the\n# document), then resolve its Paragraphs collection index so we can insert a\n# sibling list item right after it.\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute($anchorText)\nif (-not $found) {\n    throw \"Could not find the anchor paragraph: $anchorText\"\n}\n\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]0x0D, [char]0x07) -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not resolve the anchor paragraph index for: $anchorText\"\n}\n\n$anchorParagraph = $d.Paragraphs.Item($anchorIndex)\n\n# InsertParagraphAfter() splits in a new (empty) paragraph right after the\n# anchor, inheriting the anchor's paragraph formatting (the \"List Paragraph\"\n# style plus its numPr numbering), same as the sibling steps in this list.\n$anchorParagraph.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Item($anchorIndex + 1)\n$newParagraph.Range.Text = $newStepText\n"}
